$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 246, shifting existing rows 246:289 down to 247:290.
$ws.Rows.Item(246).Insert()

# Populate the newly inserted row 246 with the new weekly record.
$ws.Cells.Item(246, 1).Value = 10
$ws.Cells.Item(246, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(246, 3).Value = "La Araucanía"
$ws.Cells.Item(246, 4).Value = 44505
$ws.Cells.Item(246, 5).Value = 9
$ws.Cells.Item(246, 6).Value = 100112023
$ws.Cells.Item(246, 7).Value = "Brócoli"
$ws.Cells.Item(246, 8).Value = "Sin especificar"
$ws.Cells.Item(246, 9).Value = "Primera"
$ws.Cells.Item(246, 10).Value = 2550
$ws.Cells.Item(246, 11).Value = 800
$ws.Cells.Item(246, 12).Value = 900
$ws.Cells.Item(246, 13).Value = 849
$ws.Cells.Item(246, 14).Value = "$/unidad"
$ws.Cells.Item(246, 15).Value = "Región del Maule"
$ws.Cells.Item(246, 16).Value = 849
$ws.Cells.Item(246, 17).Value = 1
$ws.Cells.Item(246, 18).Value = "Hortaliza"
